# Daily attendance processing - 2025-10-20 13:30:43
# For every data row, the "Recorded By" column (G) lists who recorded the
# attendance. Where that list includes "System", flip the order of the
# comma-separated names so "System" (or "system") leads the list instead
# of trailing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -ne $null -and $current -like "*System*") {
        $parts = $current.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -gt 1) {
            $reversed = $parts[($parts.Length - 1)..0]
            $newValue = [string]::Join(", ", $reversed)
            $cell.Value = $newValue
        }
    }
}
